$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 309.56668
$ws.Range("I28").Value = 285.84
$ws.Range("J28").Value = 428.2
$ws.Range("K28").Value = 285.84
$ws.Range("L28").Value = 428.2
$ws.Range("M28").Value = 199.16
$ws.Range("N28").Value = -1398.2
$ws.Range("H70").Value = 4202
$ws.Range("I70").Value = 3333.3333
$ws.Range("J70").Value = 6808
$ws.Range("K70").Value = 9999.999899999999
$ws.Range("L70").Value = 20424
$ws.Range("M70").Value = -9729.999899999999
$ws.Range("N70").Value = -20964
$ws.Range("H73").Value = 4202
$ws.Range("I73").Value = 3333.3333
$ws.Range("J73").Value = 6808
$ws.Range("K73").Value = 9999.999899999999
$ws.Range("L73").Value = 20424
$ws.Range("M73").Value = -9063.999899999999
$ws.Range("N73").Value = -22296
$ws.Range("H74").Value = 6643.125
$ws.Range("I74").Value = 5189
$ws.Range("K74").Value = 5189
$ws.Range("M74").Value = -4253
$ws.Range("H77").Value = 6643.125
$ws.Range("I77").Value = 5189
$ws.Range("K77").Value = 25945
$ws.Range("M77").Value = -21265
$ws.Range("H132").Value = 14100.12
$ws.Range("I132").Value = 13458.529
$ws.Range("K132").Value = 40375.587
$ws.Range("M132").Value = -37845.587
$ws.Range("H135").Value = 1505.5416
$ws.Range("I135").Value = 1705.85
$ws.Range("J135").Value = 504
$ws.Range("K135").Value = 15352.65
$ws.Range("L135").Value = 4536
$ws.Range("M135").Value = -12817.65
$ws.Range("N135").Value = -9606

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2643816.8
$ws.Range("I132").Value = 5062.9443
$ws.Range("J132").Value = 5018695.5
$ws.Range("K132").Value = 15188.8329
$ws.Range("L132").Value = 15056086.5
$ws.Range("M132").Value = -12658.8329
$ws.Range("N132").Value = -15061146.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2523.3906
$ws.Range("I94").Value = 2768.175
$ws.Range("J94").Value = 2115.4167
$ws.Range("K94").Value = 2768.175
$ws.Range("L94").Value = 2115.4167
$ws.Range("M94").Value = -2317.175
$ws.Range("N94").Value = -3017.4167

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1396.8889
$ws.Range("I22").Value = 827.8570999999999
$ws.Range("J22").Value = 1759
$ws.Range("K22").Value = 827.8570999999999
$ws.Range("L22").Value = 1759
$ws.Range("M22").Value = -477.8570999999999
$ws.Range("N22").Value = -2459
$ws.Range("H23").Value = 5000
$ws.Range("I23").Value = 5000
$ws.Range("K23").Value = 5000
$ws.Range("M23").Value = -4760
$ws.Range("H27").Value = 5000
$ws.Range("I27").Value = 5000
$ws.Range("K27").Value = 5000
$ws.Range("M27").Value = -4808
$ws.Range("H99").Value = 6305.706
$ws.Range("I99").Value = 3143.2856
$ws.Range("J99").Value = 8519.4
$ws.Range("K99").Value = 3143.2856
$ws.Range("L99").Value = 8519.4
$ws.Range("M99").Value = -1645.2856
$ws.Range("N99").Value = -11515.4
$ws.Range("H126").Value = 6305.706
$ws.Range("I126").Value = 3143.2856
$ws.Range("J126").Value = 8519.4
$ws.Range("K126").Value = 9429.856800000001
$ws.Range("L126").Value = 25558.2
$ws.Range("M126").Value = -6959.856800000001
$ws.Range("N126").Value = -30498.2

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 261.85
$ws.Range("I2").Value = 271.8421
$ws.Range("K2").Value = 1631.0526
$ws.Range("M2").Value = -1518.0526
$ws.Range("H4").Value = 54060904
$ws.Range("I4").Value = 64134812
$ws.Range("J4").Value = 333400
$ws.Range("K4").Value = 192404436
$ws.Range("L4").Value = 1000200
$ws.Range("M4").Value = -192404324
$ws.Range("N4").Value = -1000424
$ws.Range("H131").Value = 1494.2887
$ws.Range("I131").Value = 1415
$ws.Range("J131").Value = 1495.9579
$ws.Range("K131").Value = 4245
$ws.Range("L131").Value = 4487.8737
$ws.Range("M131").Value = 795
$ws.Range("N131").Value = -14567.8737
$ws.Range("H136").Value = 2527.25
$ws.Range("I136").Value = 2316.8572
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 6950.571599999999
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -1850.571599999999
$ws.Range("N136").Value = -22200

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17321.312
$ws.Range("I70").Value = 20435.908
$ws.Range("J70").Value = 10469.2
$ws.Range("K70").Value = 20435.908
$ws.Range("L70").Value = 10469.2
$ws.Range("M70").Value = -20165.908
$ws.Range("N70").Value = -11009.2
$ws.Range("H73").Value = 17321.312
$ws.Range("I73").Value = 20435.908
$ws.Range("J73").Value = 10469.2
$ws.Range("K73").Value = 20435.908
$ws.Range("L73").Value = 10469.2
$ws.Range("M73").Value = -19499.908
$ws.Range("N73").Value = -12341.2
$ws.Range("H108").Value = 49924.75
$ws.Range("J108").Value = 49924.75
$ws.Range("L108").Value = 49924.75
$ws.Range("N108").Value = -57604.75

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3014.1904
$ws.Range("I22").Value = 2049.2856
$ws.Range("J22").Value = 4944
$ws.Range("K22").Value = 2049.2856
$ws.Range("L22").Value = 4944
$ws.Range("M22").Value = -1754.2856
$ws.Range("N22").Value = -5534
$ws.Range("H27").Value = 3014.1904
$ws.Range("I27").Value = 2049.2856
$ws.Range("J27").Value = 4944
$ws.Range("K27").Value = 2049.2856
$ws.Range("L27").Value = 4944
$ws.Range("M27").Value = -1942.2856
$ws.Range("N27").Value = -5158
$ws.Range("H34").Value = 12000
$ws.Range("J34").Value = 12000
$ws.Range("L34").Value = 12000
$ws.Range("N34").Value = -12344
$ws.Range("H46").Value = 2376.0667
$ws.Range("I46").Value = 1294.1111
$ws.Range("K46").Value = 1294.1111
$ws.Range("M46").Value = -1106.1111
$ws.Range("H93").Value = 12687.125
$ws.Range("I93").Value = 15349.4
$ws.Range("J93").Value = 8250
$ws.Range("K93").Value = 15349.4
$ws.Range("L93").Value = 8250
$ws.Range("M93").Value = -14101.4
$ws.Range("N93").Value = -10746

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("H132").Value = 5435.246
$ws.Range("I132").Value = 2942.814
$ws.Range("K132").Value = 8828.441999999999
$ws.Range("M132").Value = -6298.441999999999
$ws.Range("M34").ClearContents()
